$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.011", "39.94")
# are preserved exactly as text, matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '20.600.18'
$ws.Range("E2").Value = '  +0.87%  '

# Row 3
$ws.Range("D3").Value = '1.477.47'
$ws.Range("E3").Value = '  +0.78%  '

# Row 4
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  -0.21%  '

# Row 5
$ws.Range("D5").Value = '0.9548'
$ws.Range("E5").Value = '  +6.95%  '

# Row 6
$ws.Range("D6").Value = '280.19'
$ws.Range("E6").Value = '  -0.17%  '

# Row 7
$ws.Range("D7").Value = '0.3655'
$ws.Range("E7").Value = '  -1.63%  '

# Row 8
$ws.Range("D8").Value = '0.3065'
$ws.Range("E8").Value = '  -3.71%  '

# Row 9
$ws.Range("D9").Value = '39.94'
$ws.Range("E9").Value = '  +0.88%  '

# Row 10
$ws.Range("D10").Value = '1.062'
$ws.Range("E10").Value = '  +1.49%  '

# Row 11
$ws.Range("D11").Value = '0.06688'
$ws.Range("E11").Value = '  +1.24%  '

# Row 12
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").Value = '  -0.31%  '

# Row 13
$ws.Range("D13").Value = '5.527'
$ws.Range("E13").Value = '  -0.21%  '

# Row 14
$ws.Range("D14").Value = '18.11'
$ws.Range("E14").Value = '  +1.46%  '

# Row 15
$ws.Range("D15").Value = '6.226'
$ws.Range("E15").Value = '  +0.43%  '

# Row 16
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '0.9553'
$ws.Range("E16").Value = '  +6.53%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001036'
$ws.Range("E17").Value = '  +0.66%  '

# Row 18
$ws.Range("D18").Value = '1.476.09'
$ws.Range("E18").Value = '  +0.04%  '

# Row 19
$ws.Range("D19").Value = '0.05959'
$ws.Range("E19").Value = '  +4.95%  '

# Row 20
$ws.Range("D20").Value = '70.03'
$ws.Range("E20").Value = '  -0.50%  '

# Row 21
$ws.Range("D21").Value = '5.507'
$ws.Range("E21").Value = '  -3.00%  '

# Row 22
$ws.Range("D22").Value = '14.48'
$ws.Range("E22").Value = '  -0.55%  '

# Row 23
$ws.Range("D23").Value = '11.09'
$ws.Range("E23").Value = '  -0.77%  '

# Row 24
$ws.Range("D24").Value = '2.264'
$ws.Range("E24").Value = '  -1.31%  '

# Row 25
$ws.Range("D25").Value = '20.626.25'
$ws.Range("E25").Value = '  -0.36%  '

# Row 26
$ws.Range("D26").Value = '143.31'
$ws.Range("E26").Value = '  +4.37%  '

# Row 27
$ws.Range("D27").Value = '2.120'
$ws.Range("E27").Value = '  -6.49%  '

# Row 28
$ws.Range("E28").Value = '  -0.77%  '

# Row 29
$ws.Range("D29").Value = '1.637.32'
$ws.Range("E29").Value = '  -0.21%  '

# Row 30
$ws.Range("D30").Value = '114.08'
$ws.Range("E30").Value = '  +0.83%  '

# Row 31
$ws.Range("D31").Value = '3.968'
$ws.Range("E31").Value = '  +0.34%  '

# Row 32
$ws.Range("D32").Value = '5.021'
$ws.Range("E32").Value = '  -1.91%  '

# Row 33
$ws.Range("D33").Value = '0.8132'
$ws.Range("E33").Value = '  -2.88%  '

# Row 34
$ws.Range("D34").Value = '0.07966'
$ws.Range("E34").Value = '  +2.17%  '

# Row 35
$ws.Range("D35").Value = '1.517'
$ws.Range("E35").Value = '  +4.42%  '

# Row 36
$ws.Range("D36").Value = '1.223'
$ws.Range("E36").Value = '  +5.51%  '

# Row 37
$ws.Range("D37").Value = '0.05834'
$ws.Range("E37").Value = '  -4.29%  '

# Row 38
$ws.Range("D38").Value = '4.747'
$ws.Range("E38").Value = '  -1.85%  '

# Row 39
$ws.Range("D39").Value = '0.02053'
$ws.Range("E39").Value = '  +0.76%  '

# Row 40
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '10.41'
$ws.Range("E40").Value = '  -1.15%  '

# Row 41
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").Value = '0.9556'
$ws.Range("E41").Value = '  +4.39%  '

# Row 42
$ws.Range("D42").Value = '0.1879'
$ws.Range("E42").Value = '  +0.78%  '

# Row 43
$ws.Range("D43").Value = '7.455'
$ws.Range("E43").Value = '  +9.01%  '

# Row 44
$ws.Range("D44").Value = '0.5314'
$ws.Range("E44").Value = '  -0.69%  '

# Row 45
$ws.Range("D45").Value = '3.540'
$ws.Range("E45").Value = '  -1.09%  '

# Row 46
$ws.Range("D46").Value = '12.29'
$ws.Range("E46").Value = '  +0.20%  '

# Row 47
$ws.Range("D47").Value = '118.05'
$ws.Range("E47").Value = '  -3.48%  '

# Row 48
$ws.Range("D48").Value = '0.5196'
$ws.Range("E48").Value = '  -1.00%  '

# Row 49
$ws.Range("D49").Value = '1.825'
$ws.Range("E49").Value = '  +0.28%  '

# Row 50
$ws.Range("D50").Value = '0.06490'
$ws.Range("E50").Value = '  +1.01%  '

# Row 51
$ws.Range("D51").Value = '0.9841'
$ws.Range("E51").Value = '  -1.14%  '

